# Applies the weekly refresh of "Fecha" (D), "Volumen" (J) and price
# columns (K, L, M, P) for the Locoto / Vega Modelo de Temuco sheet.
# The underlying records are simply reshuffled across rows 2-13 (row 8
# is unchanged), so we just overwrite the cell values row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(Fecha, Volumen, PrecioMinimo, PrecioMaximo, PrecioPromedioPonderado, PrecioPorKg)
$rows = @{
    2  = @(44476, 30,  2200, 2200, 2200, 2200)
    3  = @(44484, 40,  2200, 2200, 2200, 2200)
    4  = @(44473, 140, 1600, 1600, 1600, 1600)
    5  = @(44497, 50,  2200, 2200, 2200, 2200)
    6  = @(44203, 30,  2000, 2000, 2000, 2000)
    7  = @(44483, 50,  2200, 2200, 2200, 2200)
    9  = @(44487, 50,  2200, 2200, 2200, 2200)
    10 = @(44452, 120, 2300, 2300, 2300, 2300)
    11 = @(44453, 20,  2300, 2300, 2300, 2300)
    12 = @(44474, 20,  1600, 1600, 1600, 1600)
    13 = @(44447, 75,  2200, 2200, 2200, 2200)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}
